$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (open_browser / URL): move the hyperlinked URL from E2 into D2, clear E2's text value.
$ws.Range("D2").Value = "https://opensource-demo.orangehrmlive.com/"
$ws.Range("E2").Value = ""
$ws.Hyperlinks.Add($ws.Range("D2"), "https://opensource-demo.orangehrmlive.com/") | Out-Null

# Row 3 (enter_username): id/txtUsername -> name/username
$ws.Range("C3").Value = "name"
$ws.Range("D3").Value = "username"

# Row 4 (enter_password): id/txtPassword -> name/password
$ws.Range("C4").Value = "name"
$ws.Range("D4").Value = "password"

# Row 5 (click_login): id/btnLogin -> xpath/(submit button xpath)
$ws.Range("C5").Value = "xpath"
$ws.Range("D5").Value = "//button[@type='submit' and contains(@class, 'orangehrm-login-button')]"

# Row 6 (verify_login): id/welcome/Welcome Admin -> xpath/(dashboard breadcrumb xpath)/Dashboard
$ws.Range("C6").Value = "xpath"
$ws.Range("D6").Value = "//h6[contains(@class, 'oxd-topbar-header-breadcrumb-module') and text()='Dashboard']"
$ws.Range("E6").Value = "Dashboard"

# Row heights for wrapped long-text rows
$ws.Rows("2").RowHeight = 45
$ws.Rows("5").RowHeight = 45
$ws.Rows("6").RowHeight = 45

# Column widths: D got much wider (to fit xpath text), E narrower
$ws.Columns("D").ColumnWidth = 38.140625
$ws.Columns("E").ColumnWidth = 22.42578125

# Selection moved from E11 to D11
$ws.Range("D11").Select()
